$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Calculate GPS Coordinates": add 6 new GPS readings (rows 240-245)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Calculate GPS Coordinates")

$gpsRows = @(
  @{row=240; A=234; B=29; C=52;  D=57.81;               F=90; G=25; H=55.774999999999999},
  @{row=241; A=235; B=30; C=1;   D=17.981999999999999;  F=90; G=24; H=35.771000000000001},
  @{row=242; A=236; B=30; C=2;   D=10.744999999999999;  F=90; G=25; H=55.944000000000003},
  @{row=243; A=237; B=30; C=16;  D=5.835;                F=91; G=19; H=14.887},
  @{row=244; A=238; B=30; C=9;   D=0.76700000000000002; F=91; G=19; H=40.256},
  @{row=245; A=239; B=30; C=8;   D=53.387999999999998;  F=91; G=19; H=45.902999999999999}
)

foreach ($d in $gpsRows) {
  $r = $d.row
  $ws1.Range("A$r").Value = $d.A
  $ws1.Range("B$r").Value = $d.B
  $ws1.Range("C$r").Value = $d.C
  $ws1.Range("D$r").Value = $d.D
  $ws1.Range("E$r").Formula = "=B$r+C$r/60+D$r/3600"
  $ws1.Range("E$r").HorizontalAlignment = -4152
  $ws1.Range("F$r").Value = $d.F
  $ws1.Range("G$r").Value = $d.G
  $ws1.Range("H$r").Value = $d.H
  $ws1.Range("I$r").Formula = "=-F$r-G$r/60-H$r/3600"
}

# ---------------------------------------------------------------------------
# Sheet "All Sites": add the corresponding 6 rows (237-242) and highlight
# them green, like freshly-imported site rows awaiting review.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("All Sites")

$siteRows = @(
  @{row=237; A=234; B="Absent";  C=29.882725000000001;  D=-90.432159722222224;  hasFT=$true},
  @{row=238; A=235; B="Present"; C=30.021661666666667;  D=-90.409936388888894;  hasFT=$false},
  @{row=239; A=236; B="Present"; C=30.036318055555558;  D=-90.432206666666673;  hasFT=$false},
  @{row=240; A=237; B="Absent";  C=30.2682875;           D=-91.32080194444444;  hasFT=$true},
  @{row=241; A=238; B="Absent";  C=30.150213055555554;  D=-91.32784888888888;   hasFT=$true},
  @{row=242; A=239; B="Absent";  C=30.148163333333333;  D=-91.329417499999991;  hasFT=$true}
)

$greenColor = 5287936  # RGB(0,176,80) -> FF00B050

foreach ($d in $siteRows) {
  $r = $d.row
  $ws2.Range("A$r").Value = $d.A
  $ws2.Range("B$r").Value = $d.B
  $ws2.Range("C$r").Value = $d.C
  $ws2.Range("D$r").Value = $d.D
  $ws2.Range("A$r`:D$r").Interior.Color = $greenColor
  if ($d.hasFT) {
    $ws2.Range("F$r`:T$r").Value = 0
    $ws2.Range("F$r`:T$r").Interior.Color = $greenColor
  }
}

# ---------------------------------------------------------------------------
# View state: user finished on "Calculate GPS Coordinates" selecting the
# last new cell, then switched to / ended on "All Sites" scrolled near the
# new rows with S241 selected.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("E245").Select()

$ws2.Activate()
$ws2.Range("S241").Select()
